# Applies the "assembly results correct now" update to Sheet1.
# Only raw input cells are touched; all formula cells (H/I/L/M/N columns,
# and the B38:D40 summary block) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 7 (assembly / TITAN runtime+keff block) ---
$ws.Range("B7").Value = 1.05057
$ws.Range("C7").Value = [double]"1.9000000000000001E-4"
$ws.Range("D7").Value = 1.0501199999999999
$ws.Range("E7").Value = [double]"5.0000000000000002E-5"
$ws.Range("E7").NumberFormat = "0.00E+00"
$ws.Range("F7").Value = 1.0488964300000001
$ws.Range("G7").Value = [double]"5.5286964199999999E-5"

# --- Row 15 (assembly / K20 runtime+keff cost block) ---
$ws.Range("B15").Value = 120.42
$ws.Range("D15").Value = 321.27
$ws.Range("F15").Value = 43.15

# --- Row 21 (pincell / K20 keff block) - new cells ---
$ws.Range("B21").Value = 0.27505099999999999
$ws.Range("B21").NumberFormat = "0.00E+00"
$ws.Range("C21").Value = [double]"1.8000000000000001E-4"

# --- Row 24 (assembly / K20 keff block) ---
$ws.Range("B24").Value = 1.05033
$ws.Range("D24").Value = 1.05019
$ws.Range("E24").Value = [double]"6.0000000000000002E-5"
$ws.Range("F24").Value = 1.0488605499999999
$ws.Range("G24").Value = [double]"5.5286964199999999E-5"

# --- Row 29 (pincell / cost block 2) - new cell ---
$ws.Range("B29").Value = 146.15

# --- Row 32 (assembly / cost block 2) ---
$ws.Range("B32").Value = 81.537199999999999
$ws.Range("B32").NumberFormat = "0.00E+00"
$ws.Range("D32").Value = 160.68
$ws.Range("F32").Value = 45.06

# Force a full recalculation so dependent formulas (H, I, L, M, N, and the
# B38:D40 summary rows) pick up the new inputs.
$excel.CalculateFull()

# Update the view: selection moves to C32 with no forced scroll position.
$ws.Range("C32").Select()
